# إضافة عمود جديد 'Event ' إلى Card10
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card10")

# New "Event " header in column M, matching the style of the other
# header cells (bold, bordered, centered) by copying the format of L1.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("M1").Value = "Event "

# The new column has no data yet for the existing rows - touch each
# cell so it is materialised (and the sheet's used range/dimension
# grows to include column M) while leaving it blank.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 13).Font.Bold = $false
}

# Row 8 (851-1000 tones / 10-card entry) previously had blank cells in
# F:K - fill them with the same 'nan' placeholder used by every other
# row in this sheet for missing values.
$ws.Range("F8:K8").Value = "nan"
